# Update the "想去人数" (want-to-go count) column F values across the
# four sheets to reflect the latest scrape, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        $SheetName,
        $RowValues
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Range("F$row").Value = $RowValues[$row]
    }
}

# 展览 (Exhibitions)
Set-FValues "展览" @{
    3  = 514
    4  = 487
    5  = 8900
    7  = 11445
    8  = 95
    18 = 96
    20 = 429
    22 = 1945
    23 = 745
    24 = 682
    25 = 370
    28 = 618
    29 = 1380
    34 = 1428
    36 = 317
    37 = 63
    38 = 362
    39 = 351
    40 = 42
    41 = 151
    42 = 546
    43 = 409
    44 = 130
    45 = 830
    48 = 197
    49 = 182
}

# 演出 (Performances)
Set-FValues "演出" @{
    3  = 30
    8  = 62
    23 = 54
    24 = 86
}

# 本地生活 (Local Life)
Set-FValues "本地生活" @{
    3 = 2867
    5 = 222
}

# 全部类型 (All Types)
Set-FValues "全部类型" @{
    4  = 514
    7  = 222
    8  = 487
    9  = 8900
    11 = 11445
    12 = 95
    18 = 429
    19 = 1945
    20 = 746
    21 = 683
    22 = 370
    25 = 618
    26 = 62
    28 = 1381
    34 = 1428
    36 = 317
    38 = 362
    39 = 151
    40 = 546
    41 = 409
    42 = 130
    43 = 830
    44 = 54
    48 = 197
    49 = 182
}
